$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B6").Value = "active"
$meta.Range("B8").Value = "2023-10-16T18:33:36+00:00"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")
$invariantText = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
$elements.Range("AJ1").Value = $invariantText
